# Network parameters.xlsx edit script
# 1) Row 5 (A5/B5): relabel "FRECUENCY POINTS TO BE ANALYZED" -> "ANALYSIS FREQUENCY STEP"
#    and replace the B5 value (previously "500 Hz,5 KHz") with the full comma separated
#    list of analysis frequencies from the START FRECUENCY (50 Hz) to the END FRECUENCY
#    (50 KHz) in steps of 10 Hz, formatted like "50.00 Hz" / "1.00 KHz".
# 2) Row 11 (A11:I11): clear the sample sub-network row contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Network info")

# Build the frequency list: 50 Hz -> 50000 Hz (50 KHz) step 10 Hz
$sb = New-Object System.Text.StringBuilder
$freq = 50
$first = $true
while ($freq -le 50000) {
    if ($freq -lt 1000) {
        $token = "{0:N2} Hz" -f $freq
    } else {
        $token = "{0:N2} KHz" -f ($freq / 1000)
    }
    if (-not $first) {
        [void]$sb.Append(",")
    }
    [void]$sb.Append($token)
    $first = $false
    $freq += 10
}
$freqList = $sb.ToString()

$ws.Range("A5").Value = "ANALYSIS FREQUENCY STEP"
$ws.Range("B5").Value = $freqList

# Clear the first sub-network example row (A11:I11)
$ws.Range("A11:I11").ClearContents()
